$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 532.8333
$ws.Range("I4").Value = 532.8333
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 532.8333
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -418.8333
$ws.Range("N4").ClearContents()
$ws.Range("H11").Value = 780
$ws.Range("I11").Value = 780
$ws.Range("K11").Value = 780
$ws.Range("M11").Value = -640
$ws.Range("H41").Value = 277.8
$ws.Range("I41").Value = 128.57143
$ws.Range("J41").Value = 626
$ws.Range("K41").Value = 128.57143
$ws.Range("L41").Value = 626
$ws.Range("M41").Value = 311.42857
$ws.Range("N41").Value = -1506
$ws.Range("J86").Value = 4999.8
$ws.Range("L86").Value = 4999.8
$ws.Range("N86").Value = -7245.8
$ws.Range("J89").Value = 4999.8
$ws.Range("L89").Value = 24999
$ws.Range("N89").Value = -36231
$ws.Range("H92").Value = 706.2963
$ws.Range("I92").Value = 744.8182
$ws.Range("K92").Value = 744.8182
$ws.Range("M92").Value = 503.1818
$ws.Range("H106").Value = 1777.375
$ws.Range("I106").Value = 1777.375
$ws.Range("K106").Value = 1777.375
$ws.Range("M106").Value = -1146.375
$ws.Range("H134").Value = 73211.55
$ws.Range("J134").Value = 73211.55
$ws.Range("L134").Value = 73211.55
$ws.Range("N134").Value = -83351.55
$ws.Range("H138").Value = 3246.7415
$ws.Range("I138").Value = 4659.2666
$ws.Range("J138").Value = 2754
$ws.Range("K138").Value = 13977.7998
$ws.Range("L138").Value = 8262
$ws.Range("M138").Value = -8837.799800000001
$ws.Range("N138").Value = -18542
$ws.Range("H139").Value = 69999
$ws.Range("J139").Value = 69999
$ws.Range("L139").Value = 69999
$ws.Range("N139").Value = -80279
$ws.Range("H141").Value = 3797.7058
$ws.Range("I141").Value = 1984.0667
$ws.Range("K141").Value = 5952.2001
$ws.Range("M141").Value = -772.2001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 8221.333000000001
$ws.Range("I13").Value = 8221.333000000001
$ws.Range("K13").Value = 8221.333000000001
$ws.Range("M13").Value = -8077.333000000001
$ws.Range("H32").Value = 4287.2705
$ws.Range("I32").Value = 3410.6206
$ws.Range("J32").Value = 7465.125
$ws.Range("K32").Value = 3410.6206
$ws.Range("L32").Value = 7465.125
$ws.Range("M32").Value = -3123.6206
$ws.Range("N32").Value = -8039.125
$ws.Range("H45").Value = 27180.059
$ws.Range("I45").Value = 31716.785
$ws.Range("K45").Value = 31716.785
$ws.Range("M45").Value = -31339.785
$ws.Range("H50").Value = 1020.8571
$ws.Range("I50").Value = 333
$ws.Range("J50").Value = 1536.75
$ws.Range("K50").Value = 333
$ws.Range("L50").Value = 1536.75
$ws.Range("M50").Value = 381
$ws.Range("N50").Value = -2964.75
$ws.Range("H61").Value = 2602.7568
$ws.Range("I61").Value = 2140.7144
$ws.Range("K61").Value = 2140.7144
$ws.Range("M61").Value = -1928.7144
$ws.Range("H74").Value = 156730.61
$ws.Range("I74").Value = 309948.5
$ws.Range("K74").Value = 309948.5
$ws.Range("M74").Value = -309074.5
$ws.Range("H77").Value = 156730.61
$ws.Range("I77").Value = 309948.5
$ws.Range("K77").Value = 1549742.5
$ws.Range("M77").Value = -1545374.5
$ws.Range("H136").Value = 2602.7568
$ws.Range("I136").Value = 2140.7144
$ws.Range("K136").Value = 6422.1432
$ws.Range("M136").Value = -3872.1432

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 45444.43
$ws.Range("I82").Value = 19242.375
$ws.Range("J82").Value = 80380.5
$ws.Range("K82").Value = 19242.375
$ws.Range("L82").Value = 80380.5
$ws.Range("M82").Value = -18859.375
$ws.Range("N82").Value = -81146.5
$ws.Range("H85").Value = 45444.43
$ws.Range("I85").Value = 19242.375
$ws.Range("J85").Value = 80380.5
$ws.Range("K85").Value = 19242.375
$ws.Range("L85").Value = 80380.5
$ws.Range("M85").Value = -17916.375
$ws.Range("N85").Value = -83032.5
$ws.Range("H86").Value = 2260.1904
$ws.Range("I86").Value = 2201.6875
$ws.Range("K86").Value = 2201.6875
$ws.Range("M86").Value = -1078.6875
$ws.Range("H89").Value = 2260.1904
$ws.Range("I89").Value = 2201.6875
$ws.Range("K89").Value = 11008.4375
$ws.Range("M89").Value = -5392.4375
$ws.Range("H105").Value = 7880920
$ws.Range("I105").Value = 402058.88
$ws.Range("J105").Value = 31252360
$ws.Range("K105").Value = 402058.88
$ws.Range("L105").Value = 31252360
$ws.Range("M105").Value = -400311.88
$ws.Range("N105").Value = -31255854
$ws.Range("H134").Value = 2887.256
$ws.Range("I134").Value = 2603.4849
$ws.Range("K134").Value = 7810.4547
$ws.Range("M134").Value = -5275.4547

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1499.5
$ws.Range("I16").Value = 1335.4166
$ws.Range("K16").Value = 1335.4166
$ws.Range("M16").Value = -1048.4166
$ws.Range("H31").Value = 4184.9346
$ws.Range("I31").Value = 2585.1853
$ws.Range("J31").Value = 6458.263
$ws.Range("K31").Value = 2585.1853
$ws.Range("L31").Value = 6458.263
$ws.Range("M31").Value = -2290.1853
$ws.Range("N31").Value = -7048.263
$ws.Range("H34").Value = 4184.9346
$ws.Range("I34").Value = 2585.1853
$ws.Range("J34").Value = 6458.263
$ws.Range("K34").Value = 2585.1853
$ws.Range("L34").Value = 6458.263
$ws.Range("M34").Value = -2383.1853
$ws.Range("N34").Value = -6862.263
$ws.Range("H107").Value = 985
$ws.Range("I107").Value = 631.6667
$ws.Range("K107").Value = 631.6667
$ws.Range("M107").Value = 1288.3333
$ws.Range("H113").Value = 1499.5
$ws.Range("I113").Value = 1335.4166
$ws.Range("K113").Value = 1335.4166
$ws.Range("M113").Value = 834.5834
$ws.Range("H132").Value = 2536.889
$ws.Range("I132").Value = 1868.6666
$ws.Range("K132").Value = 5605.9998
$ws.Range("M132").Value = -3075.9998
$ws.Range("H134").Value = 2978.6956
$ws.Range("I134").Value = 2595.7144
$ws.Range("K134").Value = 7787.1432
$ws.Range("M134").Value = -5252.1432

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 4198.8
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 4198.8
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 12596.4
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -12852.4
$ws.Range("H60").Value = 1114688.8
$ws.Range("I60").Value = 3334669.8
$ws.Range("K60").Value = 10004009.4
$ws.Range("M60").Value = -10003758.4
$ws.Range("H80").Value = 3423
$ws.Range("J80").Value = 3397.5
$ws.Range("L80").Value = 10192.5
$ws.Range("N80").Value = -12064.5
$ws.Range("H83").Value = 3423
$ws.Range("J83").Value = 3397.5
$ws.Range("L83").Value = 30577.5
$ws.Range("N83").Value = -39937.5
$ws.Range("H101").Value = 7750
$ws.Range("J101").Value = 7750
$ws.Range("L101").Value = 23250
$ws.Range("N101").Value = -28118
$ws.Range("H103").Value = 2182
$ws.Range("J103").Value = 2297.7
$ws.Range("L103").Value = 6893.099999999999
$ws.Range("N103").Value = -8651.099999999999
$ws.Range("H107").Value = 641.5
$ws.Range("J107").Value = 641.5
$ws.Range("L107").Value = 1924.5
$ws.Range("N107").Value = -5764.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2308.7297
$ws.Range("I132").Value = 1920.7858
$ws.Range("K132").Value = 5762.357400000001
$ws.Range("M132").Value = -3232.357400000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9999.333000000001
$ws.Range("I62").Value = 4999
$ws.Range("J62").Value = 10999.4
$ws.Range("K62").Value = 4999
$ws.Range("L62").Value = 10999.4
$ws.Range("M62").Value = -4375
$ws.Range("N62").Value = -12247.4
$ws.Range("H65").Value = 9999.333000000001
$ws.Range("I65").Value = 4999
$ws.Range("J65").Value = 10999.4
$ws.Range("K65").Value = 24995
$ws.Range("L65").Value = 54997
$ws.Range("M65").Value = -21875
$ws.Range("N65").Value = -61237
$ws.Range("H132").Value = 7330.4
$ws.Range("I132").Value = 7533.778
$ws.Range("K132").Value = 22601.334
$ws.Range("M132").Value = -20071.334
